$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The original last week "10/1-8" is split into two weeks:
#   - a new row (9) for "10/8-12"
#   - the existing row (8) relabelled to "10/1-4"
# Set the new row's label first so the shared-string table ends up with
# "10/8-12" reusing the original slot and "10/1-4" appended as new.
$ws.Range("A9").Value = "10/8-12"
$ws.Range("A8").Value = "10/1-4"

# Row 8 gains a Friday value and a Total formula.
$ws.Range("F8").Value = 0
$ws.Range("G8").Formula = "=SUM(B8:F8)"

# New row 9 only has a Monday value filled in so far.
$ws.Range("B9").Value = 2

# Selection moves to C9.
$ws.Range("C9").Select()
